$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the whole "Exercice 1 - Ligue des champions" section.
#    It spans from the paragraph right after the "sous-programmes"
#    title paragraph, through the paragraph
#    "Ecrire l'algorithme de chacun des modules envisagés."
#    (i.e. everything up to, but not including, the
#    "Exercice 2 - Distribution des cadeaux" heading).
#    The section contains a table (the TC/TM/TS example table) which
#    must be removed as well.
# -----------------------------------------------------------------

# 1a) Delete the intro paragraphs before the example table
#     ("Exercice 1 ..." through "Exemple : pour N=6 et les 3 tableaux
#     TC, TM et TS suivants :") - 20 paragraphs.
for ($i = 1; $i -le 20; $i++) {
    $d.Content.Paragraphs.Item(2).Range.Delete()
}

# 1b) Delete the example table (TC / TM / TS), which now immediately
#     follows the title paragraph.
$d.Tables.Item(1).Delete()

# 1c) Delete the remaining paragraphs of the section, from
#     "Le programme affichera :" through
#     "Ecrire l'algorithme de chacun des modules envisagés." - 8 paragraphs.
# (Use $d.Content.Paragraphs, since $d.Paragraphs can be stale right
# after a table delete.)
for ($i = 1; $i -le 8; $i++) {
    $d.Content.Paragraphs.Item(2).Range.Delete()
}

# -----------------------------------------------------------------
# 2) Consolidate a run of adjacent, identically-formatted runs into a
#    single run (purely a structural/XML cleanup; visible text is
#    unchanged).
# -----------------------------------------------------------------
$mergeRange = $d.Content
$mergeRange.Find.Execute(
    ") qui permet de retourner le résultat du calcul de l'expression ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 0
) | Out-Null

# Assigning .Text directly (instead of using Find/Replace's replacement
# text) avoids the smart-quote autocorrection and rewrites the matched
# range as a single run. Since setting identical text is a no-op, we
# first set a differing placeholder, then the final text, to force the
# run to actually be rewritten/merged.
$mergeRange.Text = "PLACEHOLDER_TEXT_FOR_MERGE"
$mergeRange.Text = ") qui permet de retourner le résultat du calcul de l'expression "
